$wb = $excel.ActiveWorkbook

$backlog = $wb.Worksheets.Item("Backlog")
$stories6 = $wb.Worksheets.Item("Stories 6")

# --- Sheet "Stories 6": insert a new task row ("Design. Seriously.") --------
# (done first so the new shared string lands ahead of the Backlog edits below,
#  matching the order new strings were appended to the workbook.)
$stories6.Rows.Item(2).Insert()
$stories6.Range("A2").Value = "Design. Seriously."
$stories6.Rows.Item(2).RowHeight = 39

# --- Sheet "Backlog": update three story/task descriptions ------------------
$backlog.Range("B121").Value = "Find equation to govern shape of a section with two different transitions"
$backlog.Range("B123").Value = "Modify vertex bender to bend models in two different parts"
$backlog.Range("B124").Value = "Track tool on a bauble will change curved track to compound curve"

# Update the selection to the new location (A4) without leaving "Stories 6"
# as the active/selected tab.
$stories6.Range("A4").Select()
$backlog.Activate()
$backlog.Range("C121").Select()
